# Update cryptos list (Price / Volume(1h) columns, and row 48/49 swap)
# per the GitHub Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.788.83'
$ws.Range("E2").Value = '  +3.63%  '
$ws.Range("D3").Value = '3.235.23'
$ws.Range("E3").Value = '  +2.22%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.52'
$ws.Range("E5").Value = '  +2.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.53'
$ws.Range("E6").Value = '  +5.53%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  -0.85%  '
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("E10").Value = '  +2.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.436'
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").Value = '3.787.48'
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.16'
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("E15").Value = '  +2.49%  '
$ws.Range("D16").Value = '60.749.75'
$ws.Range("E16").Value = '  +3.48%  '
$ws.Range("D17").Value = '3.242.63'
$ws.Range("E17").Value = '  +2.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.33'
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.34'
$ws.Range("E20").Value = '  +3.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.17'
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.528'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.01'
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("E25").Value = '  +2.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.68'
$ws.Range("E26").Value = '  +4.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").Value = '0.0₃0906'
$ws.Range("E28").Value = '  +5.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.61'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.92'
$ws.Range("E30").Value = '  +1.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.23'
$ws.Range("E31").Value = '  +3.55%  '
$ws.Range("E32").Value = '  +4.62%  '
$ws.Range("E33").Value = '  +6.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.61'
$ws.Range("E34").Value = '  +4.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.13'
$ws.Range("E35").Value = '  +0.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.43'
$ws.Range("E36").Value = '  +7.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.53'
$ws.Range("E37").Value = '  +6.26%  '
$ws.Range("D38").Value = '2.807.48'
$ws.Range("E38").Value = '  +4.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0714'
$ws.Range("E39").Value = '  +3.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0314'
$ws.Range("E40").Value = '  +8.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.72'
$ws.Range("E41").Value = '  +2.08%  '
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.97'
$ws.Range("E43").Value = '  +2.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.725'
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("E45").Value = '  +2.19%  '
$ws.Range("D46").Value = '3.273.68'
$ws.Range("E46").Value = '  +2.06%  '
$ws.Range("E47").Value = '  +2.40%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.07'
$ws.Range("E48").Value = '  +5.33%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.21'
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.805'
$ws.Range("E50").Value = '  +7.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.03%  '
